$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = [double]"0.09985113967114748"
$ws.Range("B3").Value = [double]"0.001988119311771011"
$ws.Range("C3").Value = [double]"0.0006961788887733143"
$ws.Range("D3").Value = [double]"3.02267071136162"
$ws.Range("E3").Value = [double]"0.04238283914749392"
$ws.Range("F3").Value = [double]"0.0006236289415296763"
$ws.Range("G3").Value = [double]"0.003352609682012345"
$ws.Range("H3").Value = [double]"0.1018392589829185"
$ws.Range("B4").Value = [double]"0.00764292616106737"
$ws.Range("C4").Value = [double]"0.001138267035424474"
$ws.Range("D4").Value = [double]"5.007949710942071"
$ws.Range("E4").Value = [double]"0.06671868388681791"
$ws.Range("F4").Value = [double]"0.005411956363590513"
$ws.Range("G4").Value = [double]"0.009873895958544225"
$ws.Range("H4").Value = [double]"0.1074940658322149"
$ws.Range("B5").Value = [double]"0.01148152592827088"
$ws.Range("C5").Value = [double]"0.006157564342428064"
$ws.Range("D5").Value = [double]"5.145638055242569"
$ws.Range("E5").Value = [double]"0.1565063556070586"
$ws.Range("F5").Value = [double]"-0.0005871322373460973"
$ws.Range("G5").Value = [double]"0.02355018409388785"
$ws.Range("H5").Value = [double]"0.1113326655994184"
$ws.Range("B6").Value = [double]"0.008311162695798961"
$ws.Range("C6").Value = [double]"0.003030334163803515"
$ws.Range("D6").Value = [double]"2.638315920976591"
$ws.Range("E6").Value = [double]"0.06805148546147712"
$ws.Range("F6").Value = [double]"0.002371797291782289"
$ws.Range("G6").Value = [double]"0.01425052809981564"
$ws.Range("H6").Value = [double]"0.1081623023669464"
$ws.Range("B7").Value = [double]"0.01007699447430917"
$ws.Range("C7").Value = [double]"0.003785723053045612"
$ws.Range("D7").Value = [double]"1.731141361593052"
$ws.Range("E7").Value = [double]"0.02019258107131"
$ws.Range("F7").Value = [double]"0.002657089742776302"
$ws.Range("G7").Value = [double]"0.01749689920584203"
$ws.Range("H7").Value = [double]"0.1099281341454566"
$ws.Range("B8").Value = [double]"0.01269545587341284"
$ws.Range("C8").Value = [double]"0.004709106767732963"
$ws.Range("D8").Value = [double]"2.474726883213766"
$ws.Range("E8").Value = [double]"0.06740143495831745"
$ws.Range("F8").Value = [double]"0.003465746443100233"
$ws.Range("G8").Value = [double]"0.02192516530372545"
$ws.Range("H8").Value = [double]"0.1125465955445603"
$ws.Range("B9").Value = [double]"0.02055431128896614"
$ws.Range("C9").Value = [double]"0.003352561337343317"
$ws.Range("D9").Value = [double]"4.570048414830572"
$ws.Range("E9").Value = [double]"0.06775860968567"
$ws.Range("F9").Value = [double]"0.01398339178341218"
$ws.Range("G9").Value = [double]"0.02712523079452011"
$ws.Range("H9").Value = [double]"0.1204054509601136"
$ws.Range("B10").Value = [double]"-0.09985113967114748"
$ws.Range("C10").Value = [double]"0.0005175110880237916"
$ws.Range("D10").Value = [double]"-232.5598719686458"
$ws.Range("E10").Value = [double]"0"
$ws.Range("F10").Value = [double]"-0.1008654462750203"
$ws.Range("G10").Value = [double]"-0.09883683306727471"
$ws.Range("B11").Value = [double]"-0.04531069180190018"
$ws.Range("C11").Value = [double]"0.0005653253240514305"
$ws.Range("D11").Value = [double]"-92.90461695220378"
$ws.Range("E11").Value = [double]"9.653044707023025e-173"
$ws.Range("F11").Value = [double]"-0.04641871290488364"
$ws.Range("G11").Value = [double]"-0.04420267069891675"
$ws.Range("H11").Value = [double]"0.0545404478692473"
$ws.Range("B12").Value = [double]"-0.03864793659083667"
$ws.Range("C12").Value = [double]"0.0005423361074239768"
$ws.Range("D12").Value = [double]"-84.66396795778205"
$ws.Range("E12").Value = [double]"5.897036224024574e-143"
$ws.Range("F12").Value = [double]"-0.03971089951124009"
$ws.Range("G12").Value = [double]"-0.03758497367043327"
$ws.Range("H12").Value = [double]"0.0612032030803108"
$ws.Range("B13").Value = [double]"-0.03251946179216304"
$ws.Range("C13").Value = [double]"0.0005336241642723614"
$ws.Range("D13").Value = [double]"-72.60201894693583"
$ws.Range("E13").Value = [double]"8.967103922923954e-120"
$ws.Range("F13").Value = [double]"-0.03356534956685701"
$ws.Range("G13").Value = [double]"-0.03147357401746905"
$ws.Range("H13").Value = [double]"0.06733167787898445"
$ws.Range("B14").Value = [double]"-0.02940648473962845"
$ws.Range("C14").Value = [double]"0.0005213320285015323"
$ws.Range("D14").Value = [double]"-68.60274714900342"
$ws.Range("E14").Value = [double]"2.027799064027411e-82"
$ws.Range("F14").Value = [double]"-0.03042828028660731"
$ws.Range("G14").Value = [double]"-0.02838468919264956"
$ws.Range("H14").Value = [double]"0.07044465493151902"
$ws.Range("B15").Value = [double]"-0.0249560746907821"
$ws.Range("C15").Value = [double]"0.0005184333361847483"
$ws.Range("D15").Value = [double]"-59.18894409027136"
$ws.Range("E15").Value = [double]"5.064401240803025e-16"
$ws.Range("F15").Value = [double]"-0.02597218888983772"
$ws.Range("G15").Value = [double]"-0.0239399604917265"
$ws.Range("H15").Value = [double]"0.07489506498036537"
$ws.Range("B16").Value = [double]"-0.0214606725974691"
$ws.Range("C16").Value = [double]"0.0005185489523625461"
$ws.Range("D16").Value = [double]"-50.17529295780101"
$ws.Range("E16").Value = [double]"4.446434099163968e-22"
$ws.Range("F16").Value = [double]"-0.02247701340847965"
$ws.Range("G16").Value = [double]"-0.02044433178645855"
$ws.Range("H16").Value = [double]"0.07839046707367837"
$ws.Range("B17").Value = [double]"-0.02025757036449352"
$ws.Range("C17").Value = [double]"0.0005251051937326846"
$ws.Range("D17").Value = [double]"-46.87005496496995"
$ws.Range("E17").Value = [double]"1.146486832783053e-14"
$ws.Range("F17").Value = [double]"-0.02128676121013456"
$ws.Range("G17").Value = [double]"-0.01922837951885247"
$ws.Range("H17").Value = [double]"0.07959356930665396"
$ws.Range("B18").Value = [double]"-0.01750684689290081"
$ws.Range("C18").Value = [double]"0.0005326896621379968"
$ws.Range("D18").Value = [double]"-38.76090745324535"
$ws.Range("E18").Value = [double]"4.840196289195748e-38"
$ws.Range("F18").Value = [double]"-0.01855090308023548"
$ws.Range("G18").Value = [double]"-0.01646279070556615"
$ws.Range("H18").Value = [double]"0.08234429277824666"
$ws.Range("B19").Value = [double]"-0.01469488404979446"
$ws.Range("C19").Value = [double]"0.0005281055573162014"
$ws.Range("D19").Value = [double]"-31.74063442455243"
$ws.Range("E19").Value = [double]"1.220089556426054e-21"
$ws.Range("F19").Value = [double]"-0.01572995552612867"
$ws.Range("G19").Value = [double]"-0.01365981257346024"
$ws.Range("H19").Value = [double]"0.08515625562135302"
$ws.Range("B20").Value = [double]"-0.01233485474237278"
$ws.Range("C20").Value = [double]"0.0005374723199402971"
$ws.Range("D20").Value = [double]"-25.18818088542831"
$ws.Range("E20").Value = [double]"1.30186906298213e-11"
$ws.Range("F20").Value = [double]"-0.01338828481533233"
$ws.Range("G20").Value = [double]"-0.01128142466941323"
$ws.Range("H20").Value = [double]"0.0875162849287747"
$ws.Range("B21").Value = [double]"-0.008427277330165609"
$ws.Range("C21").Value = [double]"0.0005539678871023622"
$ws.Range("D21").Value = [double]"-13.8092701281129"
$ws.Range("E21").Value = [double]"0.07609313148143519"
$ws.Range("F21").Value = [double]"-0.009513038237785987"
$ws.Range("G21").Value = [double]"-0.007341516422545229"
$ws.Range("H21").Value = [double]"0.09142386234098186"
$ws.Range("B22").Value = [double]"-0.004995674372214783"
$ws.Range("C22").Value = [double]"0.0005552342654435618"
$ws.Range("D22").Value = [double]"-6.47348131646782"
$ws.Range("E22").Value = [double]"0.0972330604518937"
$ws.Range("F22").Value = [double]"-0.006083917371593119"
$ws.Range("G22").Value = [double]"-0.003907431372836445"
$ws.Range("H22").Value = [double]"0.0948554652989327"
$ws.Range("B23").Value = [double]"-0.003187933691330872"
$ws.Range("C23").Value = [double]"0.0005674660369115787"
$ws.Range("D23").Value = [double]"-3.085064831051782"
$ws.Range("E23").Value = [double]"0.07140547236280539"
$ws.Range("F23").Value = [double]"-0.004300150597127678"
$ws.Range("G23").Value = [double]"-0.002075716785534064"
$ws.Range("H23").Value = [double]"0.0966632059798166"
$ws.Range("B24").Value = [double]"-0.00317619355577809"
$ws.Range("C24").Value = [double]"0.0005518633124045857"
$ws.Range("D24").Value = [double]"-4.190990812819697"
$ws.Range("E24").Value = [double]"0.05474402920984317"
$ws.Range("F24").Value = [double]"-0.004257829557654421"
$ws.Range("G24").Value = [double]"-0.002094557553901758"
$ws.Range("H24").Value = [double]"0.09667494611536939"
$ws.Range("B25").Value = [double]"-0.001524728413542869"
$ws.Range("C25").Value = [double]"0.0005404897300529371"
$ws.Range("D25").Value = [double]"-2.328832355709927"
$ws.Range("E25").Value = [double]"0.2014513136974785"
$ws.Range("F25").Value = [double]"-0.002584072541495681"
$ws.Range("G25").Value = [double]"-0.0004653842855900577"
$ws.Range("H25").Value = [double]"0.09832641125760461"
$ws.Range("B26").Value = [double]"0.03291407110916019"
$ws.Range("C26").Value = [double]"0.0007454848014678429"
$ws.Range("D26").Value = [double]"27.55129756683473"
$ws.Range("E26").Value = [double]"0.0143074324837879"
$ws.Range("F26").Value = [double]"0.0314529427460727"
$ws.Range("G26").Value = [double]"0.03437519947224769"
$ws.Range("H26").Value = [double]"0.1327652107803077"
